$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for the affected rows.
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = -1
$ws.Range("F12").Value = -5
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = -11
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = -2
$ws.Range("F18").Value = 1
$ws.Range("F20").Value = -1
